# Mark "command pattern" column (D) as "OK" for the rows that now fully
# implement the command pattern (no undo/redo support yet):
#   row 4  -> "2. Add hero to player"
#   row 8  -> "6. call hero skill"
#   row 10 -> "8. change name of current player"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "OK"
$ws.Range("D8").Value = "OK"
$ws.Range("D10").Value = "OK"

# Move the active selection to D9, matching the author's final cursor spot.
$ws.Range("D9").Select()
